$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.073.60'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '1.788.85'
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E4").Value = '  -0.61%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.51'
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4312'
$ws.Range("E7").Value = '  -2.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3626'
$ws.Range("E8").Value = '  -2.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.65'
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07515'
$ws.Range("E10").Value = '  -2.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.113'
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9997'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.68'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.151'
$ws.Range("E14").Value = '  -0.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.340'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").Value = '1.781.13'
$ws.Range("E16").Value = '  +1.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.97'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001064'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06332'
$ws.Range("E19").Value = '  +1.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9996'
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.950'
$ws.Range("E22").Value = '  -3.93%  '
$ws.Range("D23").Value = '28.066.32'
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.41'
$ws.Range("E24").Value = '  -2.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.162'
$ws.Range("E25").Value = '  -7.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.86'
$ws.Range("E26").Value = '  +3.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.40'
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("D28").Value = '1.987.41'
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.180'
$ws.Range("E29").Value = '  -7.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.85'
$ws.Range("E30").Value = '  -1.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.164'
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.689'
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08997'
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.513'
$ws.Range("E34").Value = '  -4.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.64'
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02322'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.088'
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6450'
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2111'
$ws.Range("E39").Value = '  -3.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06059'
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.185'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.414'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9995'
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.870'
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.56'
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5986'
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.709'
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.51'
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.988'
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06949'
$ws.Range("E51").Value = '  +0.56%  '
